# Updated symbol list (Price / Volume(1h) columns) for the GitHub Actions
# crypto-ticker refresh. Values are stored as plain text (matching the
# sheet's existing inlineStr cells), so each assignment uses a leading
# apostrophe to force text entry, then resets the cell style back to
# "Normal" to strip the quote-prefix formatting Excel applies automatically
# (keeping the cells style-less, like the original data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'274.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.85%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.56%"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'-3.08%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06315"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.18%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.934"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.52%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.377"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'44.18%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8775"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.82%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.39%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.05008"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-3.45%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07440"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.04%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02864"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-8.13%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.60%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001572"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.72%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006368"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.25%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005785"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.53%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-1.61%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.304"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.36%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-1.16%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'0.1324"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.80%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.903"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.80%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04405"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.83%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001175"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.13%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.003835"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.63%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001202"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.31%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001938"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'14.61%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.04100"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.58%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006831"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.97%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1171"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.43%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002132"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-9.52%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01149"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-8.51%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005186"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.13%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.486"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-36.91%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-11.40%"
$ws.Range("E47").Style = "Normal"
